# Update with restock suggestion
# - Forecast Comparison sheet: fill in Week_Start_Date (col B), refresh
#   Inventory Coverage (col L) and Seasonality Index (col P), replace the
#   "Sales Volume Rank" column (Q) with the recalculated "Lifecycle Stage"
#   values, and drop the now-redundant trailing "Lifecycle Stage" column (R).
#   Rows 16-17 also pick up new Stockout Risk / Reorder Urgency values.
# - Summary sheet: Max/Min Forecast Week no longer resolve to a single week,
#   so they report "N/A".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# Week start dates for W1..W16 (rows 2..17)
$weekStartDates = @(
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27",
    "2025-05-04",
    "2025-05-11",
    "2025-05-18"
)

# Inventory Coverage (column L) - counts down toward restock
$inventoryCoverage = @(14.42, 13.42, 12.42, 11.42, 10.42, 9.42, 8.42, 7.42, 6.42, 5.42, 4.42, 3.42, 2.42, 1.42, 0.42, 0)

# Seasonality Index (column P)
$seasonalityIndex = @(0.9399999999999999, 0.98, 0.9, 0.9399999999999999, 0.88, 0.9, 0.8, 1.18, 0.88, 1.19, 1.17, 0.96, 0.85, 1.02, 1.17, 0.9)

# Lifecycle Stage (replaces column Q's old "Sales Volume Rank" values)
$lifecycleStage = @("Mature","Mature","Mature","Mature","Mature","Mature","Mature","Mature","Mature","Mature","Mature","Mature","Mature","Mature","Mature","Mature")

# Column Q's header changes from "Sales Volume Rank" to "Lifecycle Stage"
# (the stat that used to live in column R).
$ws.Range("Q1").Value = "Lifecycle Stage"

# Force column B to text so the week-start dates stay as plain strings
# instead of being auto-converted to Excel date serials.
$ws.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]
    $ws.Cells.Item($row, 12).Value = $inventoryCoverage[$i]
    $ws.Cells.Item($row, 16).Value = $seasonalityIndex[$i]
    $ws.Cells.Item($row, 17).Value = $lifecycleStage[$i]
}

# Rows 16 & 17 (W15, W16) also have refreshed Stockout Risk / Reorder Urgency
$ws.Range("M16").Value = "High"
$ws.Range("N16").Value = "Urgent"
$ws.Range("M17").Value = "High"
$ws.Range("N17").Value = "Urgent"

# The old column R ("Lifecycle Stage") is now redundant since its data
# moved into column Q - delete it so the sheet ends at column Q.
$ws.Columns("R").Delete()

# Summary sheet: Max/Min Forecast Week are no longer a single deterministic
# week, so they now read "N/A".
$summary.Range("B13").Value = "N/A"
$summary.Range("B15").Value = "N/A"
